# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "69.418.21"
Set-TextCell "E2" "  -2.03%  "

Set-TextCell "D3" "3.489.09"
Set-TextCell "E3" "  -2.31%  "

Set-TextCell "D5" "611.17"
Set-TextCell "E5" "  +4.88%  "

Set-TextCell "D6" "188.34"
Set-TextCell "E6" "  -0.07%  "

Set-TextCell "D7" "0.624"
Set-TextCell "E7" "  -0.77%  "

Set-TextCell "E8" "  -0.05%  "

Set-TextCell "D9" "0.211"
Set-TextCell "E9" "  -2.81%  "

Set-TextCell "D10" "0.644"
Set-TextCell "E10" "  -1.21%  "

Set-TextCell "D11" "52.67"
Set-TextCell "E11" "  -3.74%  "

Set-TextCell "D12" "0.0000305"
Set-TextCell "E12" "  -3.77%  "

Set-TextCell "D13" "9.43"
Set-TextCell "E13" "  -1.32%  "

Set-TextCell "D14" "4.054.56"
Set-TextCell "E14" "  -2.02%  "

Set-TextCell "D15" "603.85"
Set-TextCell "E15" "  +4.05%  "

Set-TextCell "D16" "69.520.88"
Set-TextCell "E16" "  -1.93%  "

Set-TextCell "D17" "18.87"
Set-TextCell "E17" "  -1.83%  "

Set-TextCell "D18" "12.52"
Set-TextCell "E18" "  -2.05%  "

Set-TextCell "D19" "3.491.92"
Set-TextCell "E19" "  -2.30%  "

Set-TextCell "E20" "  -0.47%  "

Set-TextCell "D21" "0.980"
Set-TextCell "E21" "  -2.46%  "

Set-TextCell "D22" "17.02"
Set-TextCell "E22" "  -3.70%  "

Set-TextCell "D23" "105.83"
Set-TextCell "E23" "  +12.51%  "

Set-TextCell "D24" "4.69"
Set-TextCell "E24" "  +2.46%  "

Set-TextCell "D25" "5.11"
Set-TextCell "E25" "  +4.53%  "

Set-TextCell "D26" "3.01"
Set-TextCell "E26" "  +1.16%  "

Set-TextCell "D27" "10.90"
Set-TextCell "E27" "  -2.97%  "

Set-TextCell "D28" "9.61"
Set-TextCell "E28" "  +2.91%  "

Set-TextCell "D29" "33.25"
Set-TextCell "E29" "  +1.77%  "

Set-TextCell "D30" "6.88"
Set-TextCell "E30" "  -4.79%  "

Set-TextCell "D31" "12.46"
Set-TextCell "E31" "  +1.13%  "

Set-TextCell "D32" "4.04"
Set-TextCell "E32" "  +5.93%  "

Set-TextCell "E33" "  -2.46%  "

Set-TextCell "D34" "63.13"
Set-TextCell "E34" "  +0.06%  "

Set-TextCell "E35" "  -5.67%  "

Set-TextCell "E36" "  -0.16%  "

Set-TextCell "D37" "3.615.05"
Set-TextCell "E37" "  -0.03%  "

Set-TextCell "E38" "  +4.46%  "

Set-TextCell "D39" "0.392"
Set-TextCell "E39" "  -5.04%  "

Set-TextCell "D40" "36.53"
Set-TextCell "E40" "  -4.24%  "

Set-TextCell "B41" "Bittensor"
Set-TextCell "C41" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D41" "500.59"
Set-TextCell "E41" "  -7.16%  "

Set-TextCell "B42" "PEPE"
Set-TextCell "C42" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell "D42" "0.0₃0767"
Set-TextCell "E42" "  -6.25%  "

Set-TextCell "D43" "0.135"
Set-TextCell "E43" "  -2.98%  "

Set-TextCell "D44" "0.0459"
Set-TextCell "E44" "  -3.63%  "

Set-TextCell "D45" "2.89"
Set-TextCell "E45" "  -2.25%  "

Set-TextCell "E46" "  +1.98%  "

Set-TextCell "D47" "3.32"
Set-TextCell "E47" "  -4.82%  "

Set-TextCell "E48" "  +0.31%  "

Set-TextCell "D49" "8.72"
Set-TextCell "E49" "  -7.18%  "

Set-TextCell "D50" "130.94"
Set-TextCell "E50" "  -2.83%  "

Set-TextCell "D51" "1.34"
Set-TextCell "E51" "  -8.35%  "
